# Weekly update: insert two new price-observation rows for "Alcachofa" at
# Vega Monumental Concepción (row 30), pushing the existing rows 30-66 down
# to 32-68.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right above the current row 30 (shifts old rows
# 30..66 down to 32..68, carrying their content/formatting with them).
$ws.Rows.Item(30).Insert()
$ws.Rows.Item(30).Insert()

# --- New row 30 ---
$ws.Cells.Item(30, 1).Value = 11
$ws.Cells.Item(30, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(30, 3).Value = "Bíobío"
$ws.Cells.Item(30, 4).Value = 44803
$ws.Cells.Item(30, 5).Value = 8
$ws.Cells.Item(30, 6).Value = 100112013
$ws.Cells.Item(30, 7).Value = "Alcachofa"
$ws.Cells.Item(30, 8).Value = "Argentina(o)"
$ws.Cells.Item(30, 9).Value = "Primera"
$ws.Cells.Item(30, 10).Value = 170
$ws.Cells.Item(30, 11).Value = 13000
$ws.Cells.Item(30, 12).Value = 14000
$ws.Cells.Item(30, 13).Value = 13529
$ws.Cells.Item(30, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(30, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(30, 16).Value = 271
$ws.Cells.Item(30, 17).Value = 50
$ws.Cells.Item(30, 18).Value = "Hortaliza"

# --- New row 31 ---
$ws.Cells.Item(31, 1).Value = 11
$ws.Cells.Item(31, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(31, 3).Value = "Bíobío"
$ws.Cells.Item(31, 4).Value = 44803
$ws.Cells.Item(31, 5).Value = 8
$ws.Cells.Item(31, 6).Value = 100112013
$ws.Cells.Item(31, 7).Value = "Alcachofa"
$ws.Cells.Item(31, 8).Value = "Española"
$ws.Cells.Item(31, 9).Value = "Primera"
$ws.Cells.Item(31, 10).Value = 150
$ws.Cells.Item(31, 11).Value = 15000
$ws.Cells.Item(31, 12).Value = 16000
$ws.Cells.Item(31, 13).Value = 15467
$ws.Cells.Item(31, 14).Value = "$/caja 30 unidades"
$ws.Cells.Item(31, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(31, 16).Value = 516
$ws.Cells.Item(31, 17).Value = 30
$ws.Cells.Item(31, 18).Value = "Hortaliza"

# Make sure the D column (date) of the two new rows carries the same
# date format as the rest of the column.
$ws.Range("D30:D31").NumberFormat = $ws.Range("D32").NumberFormat()
